# Appends new daily rows (business days, skipping the 2025 US market
# holidays) to the PYPL sentiment sheet, extending the data from row 2509
# (2025-06-24) through row 2566 (2025-09-15), all with Sentiment = 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateSerials = @(45833,45834,45835,45838,45839,45840,45841,45845,45846,45847,45848,45849,45852,45853,45854,45855,45856,45859,45860,45861,45862,45863,45866,45867,45868,45869,45870,45873,45874,45875,45876,45877,45880,45881,45882,45883,45884,45887,45888,45889,45890,45891,45894,45895,45896,45897,45898,45902,45903,45904,45905,45908,45909,45910,45911,45912,45915)

# The existing date column (A) uses cell style index "2" (custom date
# number format). Grab that format from the last populated row so the new
# cells match exactly, without Excel needing to create a brand new style.
$dateFormat = $ws.Range("A2509").NumberFormat()

$startRow = 2510
for ($i = 0; $i -lt $dateSerials.Length; $i++) {
    $r = $startRow + $i
    $aCell = $ws.Cells.Item($r, 1)
    $aCell.NumberFormat = $dateFormat
    $aCell.Value = $dateSerials[$i]

    $bCell = $ws.Cells.Item($r, 2)
    $bCell.Value = 0
}
